# Reorganize the product import template header row:
#  - keep product_id / product_option_id in place (A1/B1)
#  - move name/description/meta_title/meta_description/meta_keyword to
#    immediately follow product_option_id (C1:G1)
#  - shift quantities/units/prices/image name/library/library_base_price
#    to H1:M1
#  - append the two new fields called out in the commit message
#    (percentage, multiplier) plus the rest of the new import columns
#    (attribute_ids, attribute_texts, categories, SEO) at N1:S1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "product_id"
$ws.Range("B1").Value = "product_option_id"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "description"
$ws.Range("E1").Value = "meta_title"
$ws.Range("F1").Value = "meta_description"
$ws.Range("G1").Value = "meta_keyword"
$ws.Range("H1").Value = "quantities"
$ws.Range("I1").Value = "units"
$ws.Range("J1").Value = "prices"
$ws.Range("K1").Value = "image name"
$ws.Range("L1").Value = "library (0 = no, 1 = yes)"
$ws.Range("M1").Value = "library_base_price"
$ws.Range("N1").Value = "percentage"
$ws.Range("O1").Value = "multiplier"
$ws.Range("P1").Value = "attribute_ids"
$ws.Range("Q1").Value = "attribute_texts"
$ws.Range("R1").Value = "categories"
$ws.Range("S1").Value = "SEO"

# Match the author's final selection (whole used header range selected).
$ws.Range("A1:S1").Select()

Write-Output "Header row expanded to A1:S1 with percentage/multiplier and related columns added."
